# Daily auto-push style update: insert two new hourly ranking rows at
# row 684, pushing the existing rows (old 684-725) down to 686-727, and
# populate the two freshly inserted rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 684 (shifts everything
# below down by 2, just like the authoring tool's daily append/insert).
$ws.Rows("684:685").Insert()

# Column A holds dates written as plain text ("2026/01/19", ...) rather
# than real Excel date serials, matching every other row in the sheet.
# Force Text format before assigning so Excel doesn't auto-convert the
# string into a date number.
$ws.Range("A684:A685").NumberFormat = "@"

# New row 684: 2026/01/19 (Mon), hour 23, ranking 127
$ws.Range("A684").Value = "2026/01/19"
$ws.Range("B684").Value = "月"
$ws.Range("C684").Value = 23
$ws.Range("D684").Value = 127

# New row 685: 2026/01/20 (Tue), hour 2, ranking 136
$ws.Range("A685").Value = "2026/01/20"
$ws.Range("B685").Value = "火"
$ws.Range("C685").Value = 2
$ws.Range("D685").Value = 136
